$wb = $excel.ActiveWorkbook

# Sheet ALC, row 69 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4663.091
$ws.Range("I69").Value = 4133.6665
$ws.Range("J69").Value = 5029.615
$ws.Range("K69").Value = 12400.9995
$ws.Range("L69").Value = 15088.845
$ws.Range("M69").Value = -11526.9995
$ws.Range("N69").Value = -16836.845

# Sheet ALC, row 72 (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 4663.091
$ws.Range("I72").Value = 4133.6665
$ws.Range("J72").Value = 5029.615
$ws.Range("K72").Value = 37202.9985
$ws.Range("L72").Value = 45266.535
$ws.Range("M72").Value = -32834.9985
$ws.Range("N72").Value = -54002.535

# Sheet ALC, row 103 (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1105.4166
$ws.Range("I103").Value = 821.4286
$ws.Range("J103").Value = 1503
$ws.Range("K103").Value = 2464.2858
$ws.Range("L103").Value = 4509
$ws.Range("M103").Value = -1878.2858
$ws.Range("N103").Value = -5681

# Sheet ALC, row 135 (hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1192.1936
$ws.Range("I135").Value = 1192.1936
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10729.7424
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -8194.742400000001
$ws.Range("N135").ClearContents()

# Sheet ARM, row 32 (hunk 4)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7854.4736
$ws.Range("I32").Value = 6356.6
$ws.Range("J32").Value = 10735
$ws.Range("K32").Value = 6356.6
$ws.Range("L32").Value = 10735
$ws.Range("M32").Value = -6069.6
$ws.Range("N32").Value = -11309

# Sheet ARM, row 132 (hunk 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3228.2285
$ws.Range("I132").Value = 1413.7826
$ws.Range("J132").Value = 6705.9165
$ws.Range("K132").Value = 4241.3478
$ws.Range("L132").Value = 20117.7495
$ws.Range("M132").Value = -1711.3478
$ws.Range("N132").Value = -25177.7495

# Sheet BSM, row 134 (hunk 6)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3430.1296
$ws.Range("I134").Value = 3602.25
$ws.Range("J134").Value = 2672.8
$ws.Range("K134").Value = 10806.75
$ws.Range("L134").Value = 8018.400000000001
$ws.Range("M134").Value = -8271.75
$ws.Range("N134").Value = -13088.4

# Sheet CRP, row 17 (hunk 7)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 166705840
$ws.Range("I17").Value = 15000
$ws.Range("J17").Value = 200044000
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 200044000
$ws.Range("M17").Value = -14826
$ws.Range("N17").Value = -200044348

# Sheet CRP, row 31 (hunk 8)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8733
$ws.Range("I31").Value = 2770.2727
$ws.Range("J31").Value = 13778.385
$ws.Range("K31").Value = 2770.2727
$ws.Range("L31").Value = 13778.385
$ws.Range("M31").Value = -2475.2727
$ws.Range("N31").Value = -14368.385

# Sheet CRP, row 34 (hunk 9)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8733
$ws.Range("I34").Value = 2770.2727
$ws.Range("J34").Value = 13778.385
$ws.Range("K34").Value = 2770.2727
$ws.Range("L34").Value = 13778.385
$ws.Range("M34").Value = -2568.2727
$ws.Range("N34").Value = -14182.385

# Sheet CRP, row 41 (hunk 10)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 11333.333
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 16000
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 16000
$ws.Range("M41").Value = -1572
$ws.Range("N41").Value = -16856

# Sheet CRP, row 50 (hunk 11)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 19996.25
$ws.Range("J50").Value = 19996.25
$ws.Range("L50").Value = 19996.25
$ws.Range("N50").Value = -21246.25

# Sheet CRP, row 51 (hunk 12)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 17274.75
$ws.Range("I51").Value = 12400
$ws.Range("J51").Value = 18899.666
$ws.Range("K51").Value = 12400
$ws.Range("L51").Value = 18899.666
$ws.Range("M51").Value = -11664
$ws.Range("N51").Value = -20371.666

# Sheet CRP, row 59 (hunk 13)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 24783.857
$ws.Range("J59").Value = 24783.857
$ws.Range("L59").Value = 24783.857
$ws.Range("N59").Value = -27073.857

# Sheet CRP, row 60 (hunk 14)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 23937.8
$ws.Range("J60").Value = 23937.8
$ws.Range("L60").Value = 23937.8
$ws.Range("N60").Value = -24959.8

# Sheet CRP, row 61 (hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 17274.75
$ws.Range("I61").Value = 12400
$ws.Range("J61").Value = 18899.666
$ws.Range("K61").Value = 12400
$ws.Range("L61").Value = 18899.666
$ws.Range("M61").Value = -12052
$ws.Range("N61").Value = -19595.666

# Sheet CRP, row 68 (hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 29666.334
$ws.Range("J68").Value = 29666.334
$ws.Range("L68").Value = 29666.334
$ws.Range("N68").Value = -31164.334

# Sheet CRP, row 71 (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 29666.334
$ws.Range("J71").Value = 29666.334
$ws.Range("L71").Value = 88999.00199999999
$ws.Range("N71").Value = -96487.00199999999

# Sheet CRP, row 74 (hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 21231.4
$ws.Range("J74").Value = 21231.4
$ws.Range("L74").Value = 21231.4
$ws.Range("N74").Value = -22979.4

# Sheet CRP, row 77 (hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 21231.4
$ws.Range("J77").Value = 21231.4
$ws.Range("L77").Value = 63694.2
$ws.Range("N77").Value = -72430.20000000001

# Sheet CRP, row 131 (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# Sheet CUL, row 120 (hunk 21)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 8000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 8000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 24000
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -33676

# Sheet CUL, row 121 (hunk 22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 906.6
$ws.Range("J121").Value = 1060.6571
$ws.Range("L121").Value = 3181.9713
$ws.Range("N121").Value = -5801.971299999999

# Sheet CUL, row 125 (hunk 23)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 3889.75
$ws.Range("J125").Value = 4049.0667
$ws.Range("L125").Value = 12147.2001
$ws.Range("N125").Value = -21987.2001

# Sheet CUL, row 129 (hunk 24)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2109.7144
$ws.Range("I129").Value = 1504.875
$ws.Range("J129").Value = 2916.1667
$ws.Range("K129").Value = 4514.625
$ws.Range("L129").Value = 8748.500100000001
$ws.Range("M129").Value = 485.375
$ws.Range("N129").Value = -18748.5001

# Sheet CUL, row 130 (hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 7700
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 7700
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 23100
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -33140

# Sheet CUL, row 131 (hunk 26)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1334297.2
$ws.Range("I131").Value = 6667439.5
$ws.Range("J131").Value = 1011.73334
$ws.Range("K131").Value = 20002318.5
$ws.Range("L131").Value = 3035.20002
$ws.Range("M131").Value = -19997278.5
$ws.Range("N131").Value = -13115.20002

# Sheet CUL, row 137 (hunk 27)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 21595.191
$ws.Range("I137").Value = 14212.5
$ws.Range("J137").Value = 26138.385
$ws.Range("K137").Value = 42637.5
$ws.Range("L137").Value = 78415.155
$ws.Range("M137").Value = -37537.5
$ws.Range("N137").Value = -88615.155

# Sheet GSM, row 17 (hunk 28)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 3000
$ws.Range("J17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("N17").Value = -3336

# Sheet GSM, row 63 (hunk 29)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 42110
$ws.Range("J63").Value = 42110
$ws.Range("L63").Value = 42110
$ws.Range("N63").Value = -43482

# Sheet GSM, row 66 (hunk 30)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 42110
$ws.Range("J66").Value = 42110
$ws.Range("L66").Value = 126330
$ws.Range("N66").Value = -133194

# Sheet GSM, row 70 (hunk 31)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6500.32
$ws.Range("I70").Value = 6653.0527
$ws.Range("J70").Value = 6016.6665
$ws.Range("K70").Value = 6653.0527
$ws.Range("L70").Value = 6016.6665
$ws.Range("M70").Value = -6383.0527
$ws.Range("N70").Value = -6556.6665

# Sheet GSM, row 73 (hunk 32)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6500.32
$ws.Range("I73").Value = 6653.0527
$ws.Range("J73").Value = 6016.6665
$ws.Range("K73").Value = 6653.0527
$ws.Range("L73").Value = 6016.6665
$ws.Range("M73").Value = -5717.0527
$ws.Range("N73").Value = -7888.6665

# Sheet GSM, row 126 (hunk 33)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6318.478
$ws.Range("I126").Value = 8354.933999999999
$ws.Range("J126").Value = 2500.125
$ws.Range("K126").Value = 25064.802
$ws.Range("L126").Value = 7500.375
$ws.Range("M126").Value = -22594.802
$ws.Range("N126").Value = -12440.375

# Sheet GSM, row 132 (hunk 34)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3761.4167
$ws.Range("I132").Value = 3545.4211
$ws.Range("K132").Value = 10636.2633
$ws.Range("M132").Value = -8106.263300000001

# Sheet LTW, row 21 (hunk 35)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Sheet LTW, row 136 (hunk 36)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9052.518
$ws.Range("I136").Value = 10816.5
$ws.Range("J136").Value = 7807.353
$ws.Range("K136").Value = 32449.5
$ws.Range("L136").Value = 23422.059
$ws.Range("M136").Value = -29899.5
$ws.Range("N136").Value = -28522.059

# Sheet WVR, row 132 (hunk 37)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1449
$ws.Range("I132").Value = 1010.8571
$ws.Range("K132").Value = 3032.5713
$ws.Range("M132").Value = -502.5712999999996

# Sheet WVR, row 136 (hunk 38)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 851.38464
$ws.Range("I136").Value = 536.3333
$ws.Range("J136").Value = 1901.5555
$ws.Range("K136").Value = 1608.9999
$ws.Range("L136").Value = 5704.666499999999
$ws.Range("M136").Value = 941.0001
$ws.Range("N136").Value = -10804.6665
